# Auto-generated edit script applying the scheduled-runner value refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 2
$ws.Range("J12").Value = 2
$ws.Range("L12").Value = 2
$ws.Range("N12").Value = -342

# Row 80
$ws.Range("H80").Value = 702.25
$ws.Range("I80").Value = 342
$ws.Range("K80").Value = 1026
$ws.Range("M80").Value = -28

# Row 83
$ws.Range("H83").Value = 702.25
$ws.Range("I83").Value = 342
$ws.Range("K83").Value = 3078
$ws.Range("M83").Value = 1914

# Row 88
$ws.Range("H88").Value = 1174.375
$ws.Range("I88").Value = 924.5
$ws.Range("J88").Value = 1424.25
$ws.Range("K88").Value = 924.5
$ws.Range("L88").Value = 1424.25
$ws.Range("M88").Value = -518.5
$ws.Range("N88").Value = -2236.25

# Row 91
$ws.Range("H91").Value = 1174.375
$ws.Range("I91").Value = 924.5
$ws.Range("J91").Value = 1424.25
$ws.Range("K91").Value = 924.5
$ws.Range("L91").Value = 1424.25
$ws.Range("M91").Value = 479.5
$ws.Range("N91").Value = -4232.25

# Row 111
$ws.Range("H111").Value = 3049.6667
$ws.Range("I111").Value = 2750
$ws.Range("K111").Value = 8250
$ws.Range("M111").Value = -5183

# Row 132
$ws.Range("H132").Value = 4360.5557
$ws.Range("I132").Value = 4755.625
$ws.Range("K132").Value = 14266.875
$ws.Range("M132").Value = -11736.875

# Row 135
$ws.Range("H135").Value = 3382.5557
$ws.Range("I135").Value = 3405.375
$ws.Range("J135").Value = 3200
$ws.Range("K135").Value = 30648.375
$ws.Range("L135").Value = 28800
$ws.Range("M135").Value = -28113.375
$ws.Range("N135").Value = -33870

# Row 137
$ws.Range("H137").Value = 3747.8333
$ws.Range("J137").Value = 4245
$ws.Range("L137").Value = 12735
$ws.Range("N137").Value = -17835

# Row 138
$ws.Range("H138").Value = 4192.6665
$ws.Range("I138").Value = 3394
$ws.Range("J138").Value = 5790
$ws.Range("K138").Value = 10182
$ws.Range("L138").Value = 17370
$ws.Range("M138").Value = -5042
$ws.Range("N138").Value = -27650

# Row 141
$ws.Range("H141").Value = 20558.8
$ws.Range("I141").Value = 20558.8
$ws.Range("K141").Value = 61676.39999999999
$ws.Range("M141").Value = -56496.39999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1446.5
$ws.Range("J86").Value = 1446.5
$ws.Range("L86").Value = 1446.5
$ws.Range("N86").Value = -3692.5

# Row 89
$ws.Range("H89").Value = 1446.5
$ws.Range("J89").Value = 1446.5
$ws.Range("L89").Value = 7232.5
$ws.Range("N89").Value = -18464.5

# Row 94
$ws.Range("H94").Value = 2100.4443
$ws.Range("I94").Value = 1149
$ws.Range("K94").Value = 1149
$ws.Range("M94").Value = -698

# Row 99
$ws.Range("H99").Value = 773.7778
$ws.Range("I99").Value = 653.4286
$ws.Range("K99").Value = 653.4286
$ws.Range("M99").Value = 844.5714

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 999.5
$ws.Range("J4").Value = 999.5
$ws.Range("L4").Value = 999.5
$ws.Range("N4").Value = -1223.5

# Row 107
$ws.Range("H107").Value = 632.8333
$ws.Range("I107").Value = 362.66666
$ws.Range("K107").Value = 362.66666
$ws.Range("M107").Value = 1557.33334

# Row 141
$ws.Range("H141").Value = 1090443.5
$ws.Range("J141").Value = 1090443.5
$ws.Range("L141").Value = 1090443.5
$ws.Range("N141").Value = -1100803.5

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 1360
$ws.Range("I4").Value = 1360
$ws.Range("K4").Value = 4080
$ws.Range("M4").Value = -3968

# Row 14
$ws.Range("H14").Value = 1699.5
$ws.Range("I14").Value = 1699.5
$ws.Range("K14").Value = 5098.5
$ws.Range("M14").Value = -4925.5

# Row 17
$ws.Range("H17").Value = 126.3
$ws.Range("J17").Value = 122.333336
$ws.Range("L17").Value = 367.000008
$ws.Range("N17").Value = -705.000008

# Row 48
$ws.Range("H48").Value = 4502.375
$ws.Range("I48").Value = 1669.1666
$ws.Range("K48").Value = 5007.4998
$ws.Range("M48").Value = -4757.4998

# Row 109
$ws.Range("H109").Value = 2903.4
$ws.Range("I109").Value = 2098
$ws.Range("J109").Value = 6125
$ws.Range("K109").Value = 6294
$ws.Range("L109").Value = 18375
$ws.Range("M109").Value = -5254
$ws.Range("N109").Value = -20455

$ws = $wb.Worksheets.Item("GSM")
# Row 58
$ws.Range("H58").Value = 18000
$ws.Range("I58").Value = 18000
$ws.Range("K58").Value = 18000
$ws.Range("M58").Value = -17723

# Row 70
$ws.Range("H70").Value = 8791.076999999999
$ws.Range("I70").Value = 5119.857
$ws.Range("K70").Value = 5119.857
$ws.Range("M70").Value = -4849.857

# Row 73
$ws.Range("H73").Value = 8791.076999999999
$ws.Range("I73").Value = 5119.857
$ws.Range("K73").Value = 5119.857
$ws.Range("M73").Value = -4183.857

# Row 122
$ws.Range("H122").Value = 1405.4445
$ws.Range("I122").Value = 1405.4445
$ws.Range("K122").Value = 4216.333500000001
$ws.Range("M122").Value = -1766.333500000001

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

# Row 25
$ws.Range("H25").Value = 5000000
$ws.Range("I25").Value = 5000000
$ws.Range("K25").Value = 5000000
$ws.Range("M25").Value = -4999770

# Row 39
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()

# Row 40
$ws.Range("H40").Value = 4497.3
$ws.Range("J40").Value = 9163
$ws.Range("L40").Value = 9163
$ws.Range("N40").Value = -9435

# Row 61
$ws.Range("H61").Value = 333
$ws.Range("I61").Value = 333
$ws.Range("K61").Value = 333
$ws.Range("M61").Value = -131

# Row 113
$ws.Range("H113").Value = 333
$ws.Range("I113").Value = 333
$ws.Range("K113").Value = 333
$ws.Range("M113").Value = 1837

# Row 136
$ws.Range("H136").Value = 133661.78
$ws.Range("I136").Value = 115501.164
$ws.Range("J136").Value = 169983
$ws.Range("K136").Value = 346503.492
$ws.Range("L136").Value = 509949
$ws.Range("M136").Value = -343953.492
$ws.Range("N136").Value = -515049

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -888
$ws.Range("N2").ClearContents()

# Row 100
$ws.Range("H100").Value = 400
$ws.Range("I100").Value = 400
$ws.Range("K100").Value = 800
$ws.Range("M100").Value = -259

# Row 104
$ws.Range("H104").Value = 6184
$ws.Range("J104").Value = 6184
$ws.Range("L104").Value = 6184
$ws.Range("N104").Value = -13172

# Row 132
$ws.Range("H132").Value = 1559.4
$ws.Range("I132").Value = 1441
$ws.Range("K132").Value = 4323
$ws.Range("M132").Value = -1793
